# Update the three-digit-number-divided-by-one-digit-number worksheet
# with a newly generated set of division problems.

$d = $word.ActiveDocument

$replacements = @(
    @{old = "217÷7="; new = "220÷4="},
    @{old = "236÷9="; new = "771÷3="},
    @{old = "839÷4="; new = "332÷4="},
    @{old = "838÷6="; new = "331÷5="},
    @{old = "299÷5="; new = "391÷6="},
    @{old = "632÷8="; new = "236÷4="},
    @{old = "540÷3="; new = "247÷3="},
    @{old = "103÷7="; new = "533÷6="},
    @{old = "467÷7="; new = "527÷8="},
    @{old = "170÷5="; new = "525÷9="},
    @{old = "830÷9="; new = "194÷6="},
    @{old = "845÷3="; new = "992÷9="},
    @{old = "410÷9="; new = "737÷6="},
    @{old = "575÷2="; new = "978÷3="},
    @{old = "352÷2="; new = "907÷6="},
    @{old = "753÷9="; new = "758÷4="},
    @{old = "417÷2="; new = "262÷8="},
    @{old = "511÷9="; new = "442÷7="},
    @{old = "767÷9="; new = "144÷5="},
    @{old = "318÷9="; new = "513÷4="},
    @{old = "820÷6="; new = "543÷5="},
    @{old = "803÷7="; new = "229÷6="},
    @{old = "916÷9="; new = "359÷3="},
    @{old = "322÷2="; new = "885÷6="},
    @{old = "368÷6="; new = "578÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
